$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 13, shifting existing rows 13-43 down to 14-44
$ws.Rows.Item(13).Insert()

# Copy row 14 (the just-shifted former row 13) into new row 13 as a baseline
# for all the columns that carry over unchanged (A,B,C,E,F,G,H,I,J,K,L,M,Q,R,T).
$ws.Range("A14:T14").Copy()
$ws.Range("A13:T13").PasteSpecial()

# Overwrite the columns that differ for the new row 13
$ws.Range("D13").Value = 44607
$ws.Range("N13").Value = 18000
$ws.Range("O13").Value = 19000
$ws.Range("P13").Value = 18500
$ws.Range("S13").Value = 1028
